$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.231.01'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +5.54%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.559.11'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +4.82%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '589.99'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +5.21%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '192.17'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +8.98%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.643'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.75%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.555.66'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +5.10%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.181'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.27%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.659'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.22%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '57.96'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +7.85%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000294'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +5.70%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '9.64'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +4.63%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.107.78'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +4.25%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '19.19'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.82%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.551.84'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +4.26%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.189.98'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +5.43%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.43'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +4.80%  '
$ws.Range('E20').Value = '  +0.44%  '
$ws.Range('E21').Value = '  +3.90%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '493.48'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.72%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.59'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +12.83%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '17.36'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +20.92%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '4.46'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +8.33%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '91.08'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.52%  '
$ws.Range('E27').Value = '  +4.02%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '11.15'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.50%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.28'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +6.25%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '31.91'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.69%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '7.41'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +13.04%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '12.08'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +4.84%  '
$ws.Range('E33').Value = '  +5.82%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '65.27'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +4.32%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.115'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +6.26%  '
$ws.Range('E36').Value = '  +5.09%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0₃0816'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +10.31%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.01%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '37.75'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +5.14%  '
$ws.Range('E40').Value = '  +5.39%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.58'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.60%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '3.299.96'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +6.28%  '
$ws.Range('E43').Value = '  +9.45%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.69'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +10.38%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0442'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +5.58%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.31'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.91%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.78'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +19.14%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.137'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.09%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.04'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +6.84%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.998'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.25%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.23'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +5.77%  '
